$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")
$ws.Activate()

$ws.Range("B17:AF17").Value = 0

$ws.Range("B17:AF17").Select()
